$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.751.51"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "3.780.06"
$ws.Range("E3").Value = "  -0.66%  "

Set-TextValue $ws.Range("D4") "0.995"
$ws.Range("E4").Value = "  -0.49%  "

Set-TextValue $ws.Range("D5") "599.12"
$ws.Range("E5").Value = "  +0.32%  "

Set-TextValue $ws.Range("D6") "163.13"
$ws.Range("E6").Value = "  -2.55%  "

$ws.Range("D7").Value = "3.779.16"
$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  -2.84%  "

Set-TextValue $ws.Range("D11") "0.445"
$ws.Range("E11").Value = "  -1.04%  "

Set-TextValue $ws.Range("D12") "6.65"
$ws.Range("E12").Value = "  +5.53%  "

$ws.Range("E13").Value = "  -3.46%  "

Set-TextValue $ws.Range("D14") "35.08"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").Value = "4.412.86"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "3.757.55"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("D17").Value = "67.727.34"
$ws.Range("E17").Value = "  -0.18%  "

Set-TextValue $ws.Range("D18") "18.17"
$ws.Range("E18").Value = "  -1.74%  "

$ws.Range("E19").Value = "  +1.74%  "

Set-TextValue $ws.Range("D20") "7.01"
$ws.Range("E20").Value = "  -1.00%  "

Set-TextValue $ws.Range("D21") "456.81"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("E22").Value = "  -4.30%  "

Set-TextValue $ws.Range("D23") "0.691"
$ws.Range("E23").Value = "  -1.10%  "

Set-TextValue $ws.Range("D24") "82.84"
$ws.Range("E24").Value = "  -0.58%  "

$ws.Range("E25").Value = "  -6.15%  "

Set-TextValue $ws.Range("D26") "11.83"
$ws.Range("E26").Value = "  -1.99%  "

Set-TextValue $ws.Range("D27") "2.09"
$ws.Range("E27").Value = "  -0.85%  "

Set-TextValue $ws.Range("D29") "9.84"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").Value = "3.923.17"
$ws.Range("E30").Value = "  -0.75%  "

$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "7.19"
$ws.Range("E32").Value = "  -2.55%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D33") "2.58"
$ws.Range("E33").Value = "  -6.57%  "

Set-TextValue $ws.Range("D34") "28.86"
$ws.Range("E34").Value = "  -2.15%  "

Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  -0.80%  "

Set-TextValue $ws.Range("D36") "8.95"
$ws.Range("E36").Value = "  -1.03%  "

Set-TextValue $ws.Range("D37") "0.0990"
$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("E38").Value = "  +3.85%  "

Set-TextValue $ws.Range("D39") "5.77"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  -2.25%  "

Set-TextValue $ws.Range("D41") "3.18"
$ws.Range("E41").Value = "  -6.80%  "

Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  -0.04%  "

Set-TextValue $ws.Range("D44") "43.46"
$ws.Range("E44").Value = "  +1.83%  "

Set-TextValue $ws.Range("D45") "47.22"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D46") "0.294"
$ws.Range("E46").Value = "  -1.96%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D47") "151.31"
$ws.Range("E47").Value = "  +2.33%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "8.28"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D49") "1.36"
$ws.Range("E49").Value = "  -0.93%  "

Set-TextValue $ws.Range("D50") "1.84"
$ws.Range("E50").Value = "  -0.42%  "

Set-TextValue $ws.Range("D51") "385.12"
$ws.Range("E51").Value = "  -1.77%  "
